$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Update the wording of the "material necesario" sentence.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Se verifico que se tenía el material necesario para llevar acabo los servicios a brindar.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Se verifico la disponibilidad de recursos y que no se cuentan problemas con los mismos. ",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2. The two blank paragraphs that used to sit right after that
#    sentence are removed (they get relocated to the end of the
#    section, right after the paragraph holding the _GoBack bookmark).
# ---------------------------------------------------------------------
$targetIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Se verifico la disponibilidad de recursos*") {
        $targetIdx = $i
        break
    }
}

$firstBlank = $d.Paragraphs.Item($targetIdx + 1)
$delStart = $firstBlank.Range.Start
$delEnd = $firstBlank.Range.End

# Deleting the same (1-char) paragraph-mark range twice collapses the
# two following blank paragraphs into the text paragraph above them.
$d.Range($delStart, $delEnd).Delete()
$d.Range($delStart, $delEnd).Delete()

# ---------------------------------------------------------------------
# 3. Two new blank paragraphs (identical formatting to the ones that
#    were removed) are appended right after the paragraph that holds
#    the _GoBack bookmark.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$anchorStart = $bm.Start

$anchorIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $anchorStart) {
        $anchorIdx = $i
        break
    }
}

$insertAfterIdx = $anchorIdx
for ($n = 0; $n -lt 2; $n++) {
    $anchorPara = $d.Paragraphs.Item($insertAfterIdx)
    $anchorPara.Range.InsertParagraphAfter()

    # The freshly-created paragraph carries over an (empty) run with the
    # neighbouring rPr; type a throw-away character and remove it again
    # so the saved XML ends up with a bare <w:p><w:pPr>...</w:pPr></w:p>,
    # matching the sibling blank paragraphs.
    $newIdx = $insertAfterIdx + 1
    $newPara = $d.Paragraphs.Item($newIdx)
    $r = $newPara.Range
    $r.End = $r.End - 1
    $r.Text = "X"

    $newPara2 = $d.Paragraphs.Item($newIdx)
    $rClear = $d.Range($newPara2.Range.Start, $newPara2.Range.Start + 1)
    $rClear.Text = ""

    $insertAfterIdx = $newIdx
}
